# Refactor build and phylogeny
# Populate previously-empty host_species (column D) cells for several
# cell-culture-derived / host-associated AAV reference sequences.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("REFSET")

$ws.Range("D6").Value = "Cell culture"
$ws.Range("D7").Value = "Cell culture"
$ws.Range("D8").Value = "Cell culture"
$ws.Range("D10").Value = "Myotis ricketti"
$ws.Range("D9").Value = "Zalophus californianus"
$ws.Range("D29").Value = "Bos taurus"

# Update selection to match the saved view state
$ws.Range("F18").Select()
